# Updates the cryptos worksheet with freshly scraped price/volume data.
# Equivalent to the automated "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing text storage (the source sheet
# stores every Price/Volume figure as text, even when it is numeric-looking), then
# restore the cell's original (default) style so no formatting changes leak in.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 47 and 48 swapped position in the source ranking (Stacks moved above Arweave),
# so update name/link columns for those two rows explicitly.
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"

$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"

# New Price (column D) / Volume(1h) (column E) values per row.
$ws.Range("D2").Value = "68.415.49"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "3.751.69"
$ws.Range("E3").Value = "  -0.57%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.11%  "
Set-TextValue "D5" "595.64"
$ws.Range("E5").Value = "  -0.14%  "
Set-TextValue "D6" "166.70"
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("D7").Value = "3.746.03"
$ws.Range("E7").Value = "  -0.67%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("E10").Value = "  -3.02%  "
Set-TextValue "D11" "6.48"
$ws.Range("E11").Value = "  +0.34%  "
Set-TextValue "D12" "0.447"
$ws.Range("E12").Value = "  -1.17%  "
$ws.Range("E13").Value = "  -6.30%  "
Set-TextValue "D14" "36.03"
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").Value = "4.380.44"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").Value = "3.727.74"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").Value = "68.412.13"
$ws.Range("E17").Value = "  +1.13%  "
Set-TextValue "D18" "17.83"
$ws.Range("E18").Value = "  -4.09%  "
$ws.Range("E19").Value = "  -2.54%  "
$ws.Range("E20").Value = "  -0.19%  "
Set-TextValue "D21" "10.69"
$ws.Range("E21").Value = "  +1.06%  "
Set-TextValue "D22" "468.14"
$ws.Range("E22").Value = "  +0.07%  "
Set-TextValue "D23" "0.697"
$ws.Range("E23").Value = "  -2.86%  "
Set-TextValue "D24" "84.45"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("E25").Value = "  -2.29%  "
$ws.Range("E26").Value = "  -0.59%  "
Set-TextValue "D27" "12.01"
$ws.Range("E27").Value = "  -1.42%  "
Set-TextValue "D28" "10.08"
$ws.Range("E28").Value = "  -1.98%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "3.899.05"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("E31").Value = "  -4.66%  "
Set-TextValue "D32" "7.28"
$ws.Range("E32").Value = "  -4.60%  "
Set-TextValue "D33" "29.84"
$ws.Range("E33").Value = "  -2.10%  "
$ws.Range("E34").Value = "  -2.20%  "
Set-TextValue "D35" "9.20"
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("D37").Value = "3.707.50"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("E38").Value = "  -2.75%  "
Set-TextValue "D39" "3.39"
$ws.Range("E39").Value = "  -11.37%  "
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("E41").Value = "  -0.44%  "
Set-TextValue "D42" "5.80"
$ws.Range("E42").Value = "  +0.18%  "
Set-TextValue "D43" "0.999"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("E45").Value = "  -2.19%  "
Set-TextValue "D46" "8.60"
$ws.Range("E46").Value = "  -0.88%  "
Set-TextValue "D47" "1.93"
$ws.Range("E47").Value = "  -0.78%  "
Set-TextValue "D48" "42.95"
$ws.Range("E48").Value = "  +9.90%  "
Set-TextValue "D49" "45.77"
$ws.Range("E49").Value = "  -0.02%  "
Set-TextValue "D50" "145.70"
$ws.Range("E50").Value = "  +4.10%  "
Set-TextValue "D51" "390.35"
$ws.Range("E51").Value = "  -1.28%  "
